# Recreate the stub "tuple row" formulas from column H (using C/D/E/F) into a
# new column P (using K/L/M/N) so the second table mirrors the first table's
# "meaningful data type" stub output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column P: row 9 gets its own (non-shared) formula, rows 10:17 are
#     filled from row 10's formula (creates the shared formula group, same
#     pattern Excel itself used for H10:H17 / O10:O16). ---
$ws.Range("P9").Formula  = '="("&"''"&K9&"''"&", "&"''"&L9&"''"&", "&M9&", "&N9&")"&",\"'
$ws.Range("P10:P17").Formula = '="("&"''"&K10&"''"&", "&"''"&L10&"''"&", "&M10&", "&N10&")"&",\"'

# --- Column I is narrower now ---
$ws.Columns.Item(9).ColumnWidth = 3.65

# --- Selection moved ---
$ws.Range("J22").Select()
